$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 20 de Agosto de 2020 a las 13:13"

# Update statistics for countries whose figures changed (country stays in place)
# Row 4 - Estados Unidos
$ws.Range("B4").Value = 5701390
$ws.Range("C4").Value = 459
$ws.Range("E4").Value = 2461812
$ws.Range("G4").Value = 28
$ws.Range("H4").Value = 176365

# Row 14 - Iran
$ws.Range("B14").Value = 352558
$ws.Range("C14").Value = 2279
$ws.Range("D14").Value = 304236
$ws.Range("E14").Value = 28058
$ws.Range("G14").Value = 139
$ws.Range("H14").Value = 20264

# Row 19 - Banglades
$ws.Range("B19").Value = 287959
$ws.Range("C19").Value = 2868
$ws.Range("D19").Value = 168991
$ws.Range("E19").Value = 115146
$ws.Range("G19").Value = 41
$ws.Range("H19").Value = 3822

# Row 42 - Rumania
$ws.Range("B42").Value = 74963
$ws.Range("C42").Value = 1346
$ws.Range("D42").Value = 34196
$ws.Range("E42").Value = 37613
$ws.Range("G42").Value = 48
$ws.Range("H42").Value = 3154

# Row 59 - Suiza
$ws.Range("B59").Value = 39026
$ws.Range("C59").Value = 266
$ws.Range("E59").Value = 3228
$ws.Range("G59").Value = 2
$ws.Range("H59").Value = 1998

# Row 69 - Nepal
$ws.Range("B69").Value = 29645
$ws.Range("C69").Value = 707
$ws.Range("D69").Value = 17964
$ws.Range("E69").Value = 11555
$ws.Range("G69").Value = 6
$ws.Range("H69").Value = 126

# Row 73 - El Salvador
$ws.Range("B73").Value = 23964
$ws.Range("C73").Value = 247
$ws.Range("D73").Value = 11558
$ws.Range("E73").Value = 11766

# Row 76 - Estado de Palestina
$ws.Range("B76").Value = 17989
$ws.Range("C76").Value = 383
$ws.Range("D76").Value = 10682
$ws.Range("E76").Value = 7188

# Row 82 - Madagascar
$ws.Range("B82").Value = 14154
$ws.Range("C82").Value = 80
$ws.Range("D82").Value = 13038
$ws.Range("E82").Value = 939
$ws.Range("G82").Value = 4
$ws.Range("H82").Value = 177

# Row 85 - Senegal
$ws.Range("B85").Value = 12559
$ws.Range("C85").Value = 113
$ws.Range("D85").Value = 8050
$ws.Range("E85").Value = 4248
$ws.Range("G85").Value = 3
$ws.Range("H85").Value = 261

# Row 90 - Consejo Danes para los Refugiados
$ws.Range("B90").Value = 9757
$ws.Range("C90").Value = 16
$ws.Range("E90").Value = 615
$ws.Range("G90").Value = 1
$ws.Range("H90").Value = 247

# Row 112 - Hong Kong
$ws.Range("B112").Value = 4605
$ws.Range("C112").Value = 18
$ws.Range("D112").Value = 3827
$ws.Range("E112").Value = 705
$ws.Range("G112").Value = 1
$ws.Range("H112").Value = 73

# Row 132 - Gambia
$ws.Range("B132").Value = 2401
$ws.Range("C132").Value = 113
$ws.Range("E132").Value = 1885
$ws.Range("G132").Value = 4
$ws.Range("H132").Value = 81

# Row 160 - Lesoto
$ws.Range("D160").Value = 423
$ws.Range("E160").Value = 493

# Rows 144-146: Malta's updated figures push it above Uruguay and Jordania
# in the descending sort by "Casos totales", so the three rows shift.
$ws.Range("A144").Value = "Malta"
$ws.Range("B144").Value = 1510
$ws.Range("C144").Value = 40
$ws.Range("D144").Value = 802
$ws.Range("E144").Value = 699
$ws.Range("H144").Value = 9

$ws.Range("A145").Value = "Uruguay"
$ws.Range("B145").Value = 1493
$ws.Range("D145").Value = 1228
$ws.Range("E145").Value = 225
$ws.Range("H145").Value = 40

$ws.Range("A146").Value = "Jordania"
$ws.Range("B146").Value = 1482
$ws.Range("D146").Value = 1259
$ws.Range("E146").Value = 212
$ws.Range("H146").Value = 11

# Rows 213-214: Montserrat's updated figures push it above Islas Malvinas
$ws.Range("A213").Value = "Montserrat"
$ws.Range("D213").Value = 12
$ws.Range("H213").Value = 1

$ws.Range("A214").Value = "Islas Malvinas"
$ws.Range("D214").Value = 13
$ws.Range("H214").Value = 0
